$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: force Text number-format on cells whose new value would otherwise
# be auto-detected as a number by Excel, so the literal text is preserved exactly.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: write the new cell values/text
$ws.Range("D2").Value = "26.187.90"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").Value = "1.657.08"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "216.50"
$ws.Range("E5").Value = "  -3.42%  "
$ws.Range("D6").Value = "0.5142"
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").Value = "0.06431"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "1.660.43"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "4.302"
$ws.Range("E13").Value = "  -4.68%  "
$ws.Range("D14").Value = "1.884.52"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").Value = "0.5545"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "0.0₅8039"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "64.27"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "26.206.79"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "211.38"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "4.416"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  -3.38%  "
$ws.Range("D23").Value = "6.017"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "144.29"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "1.724"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("D28").Value = "6.980"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "15.78"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "0.05137"
$ws.Range("E30").Value = "  -4.45%  "
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("D32").Value = "3.351"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").Value = "3.232"
$ws.Range("E33").Value = "  -5.84%  "
$ws.Range("D34").Value = "1.568"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "2.755"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.9304"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.372"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "0.5713"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").Value = "1.164.92"
$ws.Range("E39").Value = "  +10.93%  "
$ws.Range("D40").Value = "0.01595"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "0.8392"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "5.664"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.795.39"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "55.85"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "7.879"
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").Value = "0.05059"
$ws.Range("E51").Value = "  -3.36%  "

# Step 3: clear the temporary Text formatting so cell styling matches the original
# (no explicit style on these data cells), while values remain stored as text.
$ws.Range("D5,D6,D7,D8,D9,D11,D13,D15,D17,D20,D21,D22,D23,D24,D25,D26,D28,D29,D30,D32,D33,D34,D35,D36,D37,D38,D40,D42,D43,D48,D49,D50,D51").ClearFormats()
